$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (Week 7 header) - add total time in B22, matching the bold style of A22
$ws.Range("B22").Value = "4 uur"
$ws.Range("B22").Font.Bold = $true

# Row 24 (new "Week 9" header row), bold like other week header cells (e.g. A22)
$ws.Range("A24").Value = "Week 9"
$ws.Range("A24").Font.Bold = $true

# Row 25 (new log entry row)
$ws.Range("A25").Value = 43567
$ws.Range("A25").NumberFormat = $ws.Range("A23").NumberFormat
$ws.Range("B25").Value = "1 uur 5 minuten"
$ws.Range("C25").Value = "cloud firestore toevoegen aan proef project, data versturen naar firestore, data ophalen uit datastore"

# Row 27 (new row with hyperlink in Q27)
$ws.Range("Q27").Value = "https://firebase.google.com/docs/firestore/quickstart?authuser=0"
$ws.Hyperlinks.Add($ws.Range("Q27"), "https://firebase.google.com/docs/firestore/quickstart?authuser=0")
$ws.Range("Q27").Style = $ws.Range("Q26").Style

# Update selection to match the final edit position
$ws.Range("Q27").Select() | Out-Null
